# working_hours.xlsx: append four new time entries (rows 153-156), pushing
# the previous blank-separator + summary rows (formerly 153-156) down to
# 157-160, and update the summary formulas to cover the extended range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move the blank separator + 3 summary rows down by 4 rows ---------
# Insert 4 blank rows right before the old separator row (153) so the
# existing separator/summary rows shift from 153-156 to 157-160.
$ws.Range("A153:H156").Insert() | Out-Null

# --- 2. Fill in the 4 new data rows (153-156) -----------------------------
$ws.Range("A153").Value = 2014
$ws.Range("B153").Value = 7
$ws.Range("C153").Value = 25
$ws.Range("D153").Value = 10 / 24
$ws.Range("E153").Value = 11 / 24
$ws.Range("F153").Formula = "=(E153-D153)*24*60"
$ws.Range("G153").Formula = "=F153/60"

$ws.Range("A154").Value = 2014
$ws.Range("B154").Value = 7
$ws.Range("C154").Value = 25
$ws.Range("D154").Value = 16.75 / 24
$ws.Range("E154").Value = 17.25 / 24
$ws.Range("F154").Formula = "=(E154-D154)*24*60"
$ws.Range("G154").Formula = "=F154/60"

$ws.Range("A155").Value = 2014
$ws.Range("B155").Value = 7
$ws.Range("C155").Value = 26
$ws.Range("D155").Value = 9 / 24
$ws.Range("E155").Value = 10 / 24
$ws.Range("F155").Formula = "=(E155-D155)*24*60"
$ws.Range("G155").Formula = "=F155/60"

$ws.Range("A156").Value = 2014
$ws.Range("B156").Value = 7
$ws.Range("C156").Value = 26
$ws.Range("D156").Value = 11 / 24
$ws.Range("E156").Value = 12 / 24
$ws.Range("F156").Formula = "=(E156-D156)*24*60"
$ws.Range("G156").Formula = "=F156/60"

# --- 3. Update the summary formulas (now on rows 158-160) to cover the ---
# --- extended data range (F2:F156) ---------------------------------------
$ws.Range("F158").Formula = "=SUM(F2:F156)"
$ws.Range("F159").Formula = "=F158/60"
$ws.Range("F160").Formula = "=F159/38.5"

# --- 4. Update the selected cell shown in the sheet view -----------------
$ws.Range("F156").Select() | Out-Null

$wb.Save()
